$d = $word.ActiveDocument

# 1) "(34)" -> "(45)" — the emergency-light count.
$rng = $d.Content
$rng.Find.Execute("certifica que las (", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.End
$rngNum = $d.Range($start, $start + 2)
$rngNum.Text = "45"

# 2) Insert the facility name "PVEA CUZCO2" where it was blank, before the
#    comma that follows "posee la".
$rng = $d.Content
$rng.Find.Execute("posee la ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.End
$rngName = $d.Range($start, $start)
$rngName.Text = "PVEA CUZCO2"
$rngName.Bold = 1
$rngName.Font.Bold = 1

# Also bold the comma right after the inserted name.
$rng = $d.Content
$rng.Find.Execute("PVEA CUZCO2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.End
$rngComma = $d.Range($start, $start + 1)
$rngComma.Bold = 1
$rngComma.Font.Bold = 1

# 3) Bold the address "Av. Nicolás Ayllón 836" and the comma after it.
$rng = $d.Content
$rng.Find.Execute("Av. Nicolás Ayllón 836", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Bold = 1
$rng.Font.Bold = 1
$start = $rng.End
$rngComma2 = $d.Range($start, $start + 1)
$rngComma2.Bold = 1
$rngComma2.Font.Bold = 1

# 4) District / province / department block.
$d.Content.Find.Execute("en el distrito de WANCHAQ, provincia de ESPINAR y departamento de CUSCO,", $true, $false, $false, $false, $false, $true, 1, $false, "en el distrito de HUACHON, provincia de PASCO y departamento de PASCO,", 2)

# 5) Certificate date.
$d.Content.Find.Execute("20 de enero de 2026", $true, $false, $false, $false, $false, $true, 1, $false, "21 de enero de 2026", 2)
